# Daniel's assessment merged into the workbook structure.
# Fills in the "Self assesment" rows for both criterion blocks on the
# "Peer  and self assessment" sheet with grade + example text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Peer  and self assessment")

# Criterion 1 Online collaboration -> Self assesment row (row 2)
$ws.Range("B2").Value = "Insufficient"
$ws.Range("C2").Value = "Long reponse time 48H+ in some cases. Did not participate in first meeting."

# Criterion 1 International Collaboration -> Self assesment row (row 15)
$ws.Range("B15").Value = "Good"
$ws.Range("C15").Value = "Active collaborator, motivated"

# Move the active selection to reflect where the editor ended up working
[void]$ws.Activate()
[void]$excel.Goto($ws.Range("B13"), $true)
[void]$ws.Range("C15").Select()
